# "multi-ref started to work" - update MTTR (column J) for each machine row
# from 1200 to 3600 on the "Line Data" sheet, and refresh the row heights
# (rows 1-27) to their new auto-fit values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line Data")

# Update MTTR values (column J) for every machine data row (2-16): 1200 -> 3600
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 3600
}

# Refresh row heights to match the post-edit layout
$ws.Rows.Item(1).RowHeight = 18

for ($r = 2; $r -le 15; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}

$ws.Rows.Item(16).RowHeight = 19.5

for ($r = 17; $r -le 27; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
